# Update Model Component class diagram
#
# The AddressBook sample class diagram is renamed to a TaskManager
# sample: the shape labels in the UML-ish diagram on slide 1 are
# updated to reflect the new domain names.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shapes are looked up by their stable cNvPr Id (not by positional
# index) so the replacements are robust to any shape re-ordering.
function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# Simple single-run shape labels.
$renames = @{
    46 = "TaskManager"       # was AddressBook
    49 = "UniqueTaskList"    # was UniquePersonList
    62 = "Task"              # was Person
    80 = "Note"              # was Phone
    83 = "Priority"          # was Email
    85 = "Status"            # was Address
}

foreach ($id in $renames.Keys) {
    $shp = Get-ShapeById $s $id
    $shp.TextFrame.TextRange.Text = $renames[$id]
}

# Two-run shape labels: run 1 is the "<<interface>>" caption, run 2 is
# the class name after the line break - only the name run changes.
$runRenames = @{
    72  = "ReadOnlyTask"        # was ReadOnlyPerson
    100 = "ReadOnlyTaskManager" # was ReadOnlyAddressBook
}

foreach ($id in $runRenames.Keys) {
    $shp = Get-ShapeById $s $id
    $shp.TextFrame.TextRange.Runs(2).Text = $runRenames[$id]
}

# Slide master background: switch from the themed background reference
# to an explicit solid fill.
$bgFill = $p.SlideMaster.Background.Fill
$bgFill.Solid()
$bgFill.ForeColor.RGB = 16777215
